# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect newly generated output numbers.

$wb = $excel.ActiveWorkbook

# Mapping of row number -> new value for column F
$updates = @{
    2  = 373
    4  = 10837
    6  = 980
    7  = 169
    8  = 1341
    9  = 8301
    10 = 41
    12 = 588
    15 = 3306
    17 = 328
    18 = 33
    19 = 793
    22 = 287
    23 = 115
    24 = 1789
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
